# Apply the "Updated cryptos list" data refresh to the crypto price table.
# For each changed row: update Coin/Link (only for the two swapped rows),
# Price (D) and Volume(1h) (E) to the new scraped values.
#
# Price values are text (not numbers) in this sheet. Several new Price
# strings look like plain numbers/decimals (e.g. "4.106", "0.9993"), and
# Excel would silently coerce a bare numeric-looking string into a real
# number. To keep them as text (matching the original formatting, e.g.
# trailing zeros like "1.000" or "83.00"), such values are entered with a
# leading apostrophe (Excel's literal "treat as text" prefix) and the
# cell style is immediately reset to Normal so no extra formatting sticks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.419.32'
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").Value = '1.837.60'
$ws.Range("E3").Value = '  -0.71%  '

$ws.Range("D4").Value = '''0.9993'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.52%  '

$ws.Range("D5").Value = '''243.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.63%  '

$ws.Range("D6").Value = '''0.6250'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.89%  '

$ws.Range("D7").Value = '''1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.41%  '

$ws.Range("D8").Value = '''0.07394'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.84%  '

$ws.Range("D9").Value = '''0.2927'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("E10").Value = '  +0.90%  '

$ws.Range("E11").Value = '  -1.45%  '

$ws.Range("D12").Value = '1.836.69'
$ws.Range("E12").Value = '  -0.44%  '

$ws.Range("D13").Value = '''5.008'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.30%  '

$ws.Range("D14").Value = '''0.6741'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.16%  '

$ws.Range("D15").Value = '''83.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.40%  '

$ws.Range("D16").Value = '''0.000009318'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.86%  '

$ws.Range("D17").Value = '''5.877'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("D18").Value = '29.393.23'
$ws.Range("E18").Value = '  +0.16%  '

$ws.Range("D19").Value = '2.086.08'
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").Value = '''238.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").Value = '''12.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.67%  '

$ws.Range("D22").Value = '''1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.49%  '

$ws.Range("E23").Value = '  +2.11%  '

$ws.Range("D24").Value = '''1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.65%  '

$ws.Range("D25").Value = '''158.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.18%  '

$ws.Range("D26").Value = '''0.1414'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.27%  '

$ws.Range("D27").Value = '''8.477'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.74%  '

$ws.Range("D28").Value = '''17.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.16%  '

$ws.Range("D29").Value = '''0.06142'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.42%  '

$ws.Range("D30").Value = '''1.496'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.89%  '

$ws.Range("D31").Value = '''1.239'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.37%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''4.106'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.57%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '''4.086'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.93%  '

$ws.Range("E34").Value = '  +0.57%  '

$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").Value = '''0.7253'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.45%  '

$ws.Range("D37").Value = '''2.608'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.47%  '

$ws.Range("D38").Value = '''2.886'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.61%  '

$ws.Range("D39").Value = '1.216.63'
$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").Value = '''0.01761'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.52%  '

$ws.Range("D41").Value = '''6.295'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.34%  '

$ws.Range("D42").Value = '''0.9134'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.16%  '

$ws.Range("D43").Value = '''1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.30%  '

$ws.Range("D44").Value = '2.001.08'
$ws.Range("E44").Value = '  +0.71%  '

$ws.Range("D45").Value = '''101.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("D46").Value = '''65.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.35%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.5067'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.73%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '''0.00000000119'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.39%  '

$ws.Range("D49").Value = '''9.223'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.55%  '

$ws.Range("D50").Value = '''0.4050'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("D51").Value = '''0.1138'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.22%  '
